# Auto-generated Excel COM-interop script to apply the diff
# (update F/G values for rows 272-442, and append new rows 443-445)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated values in existing rows (F/G columns) ---
$ws.Range("F272").Value = 30228
$ws.Range("G272").Value = 1563
$ws.Range("F273").Value = 32048
$ws.Range("G273").Value = 1650
$ws.Range("F274").Value = 28719
$ws.Range("F275").Value = 30519
$ws.Range("F276").Value = 11501
$ws.Range("F277").Value = 3418
$ws.Range("F278").Value = 31079
$ws.Range("F279").Value = 42696
$ws.Range("F280").Value = 34883
$ws.Range("F281").Value = 46820
$ws.Range("F282").Value = 46430
$ws.Range("F285").Value = 42762
$ws.Range("F286").Value = 55327
$ws.Range("F287").Value = 59100
$ws.Range("F288").Value = 59316
$ws.Range("F289").Value = 63078
$ws.Range("F290").Value = 17757
$ws.Range("F291").Value = 15271
$ws.Range("F292").Value = 83152
$ws.Range("F293").Value = 83383
$ws.Range("F294").Value = 94052
$ws.Range("F295").Value = 17345
$ws.Range("F296").Value = 1888
$ws.Range("F297").Value = 2445
$ws.Range("F298").Value = 3286
$ws.Range("F299").Value = 66439
$ws.Range("F300").Value = 73082
$ws.Range("F301").Value = 72404
$ws.Range("F302").Value = 78934
$ws.Range("F303").Value = 9728
$ws.Range("F304").Value = 6178
$ws.Range("F307").Value = 75631
$ws.Range("F308").Value = 15806
$ws.Range("F309").Value = 78717
$ws.Range("F310").Value = 79402
$ws.Range("F312").Value = 28329
$ws.Range("F313").Value = 76621
$ws.Range("F315").Value = 56884
$ws.Range("F316").Value = 50906
$ws.Range("F317").Value = 63890
$ws.Range("F319").Value = 41414
$ws.Range("F320").Value = 74000
$ws.Range("F321").Value = 90962
$ws.Range("F322").Value = 110157
$ws.Range("F323").Value = 217657
$ws.Range("F324").Value = 250187
$ws.Range("F326").Value = 418289
$ws.Range("F327").Value = 224318
$ws.Range("F328").Value = 181255
$ws.Range("F329").Value = 73499
$ws.Range("F330").Value = 71612
$ws.Range("F334").Value = 193326
$ws.Range("F335").Value = 150551
$ws.Range("F336").Value = 81994
$ws.Range("F338").Value = 221438
$ws.Range("F339").Value = 662697
$ws.Range("F340").Value = 387256
$ws.Range("F341").Value = 283506
$ws.Range("F342").Value = 178902
$ws.Range("F343").Value = 134112
$ws.Range("F344").Value = 136063
$ws.Range("G344").Value = 2486
$ws.Range("F345").Value = 292313
$ws.Range("F346").Value = 675026
$ws.Range("F347").Value = 346654
$ws.Range("F348").Value = 232919
$ws.Range("F349").Value = 159879
$ws.Range("F350").Value = 127600
$ws.Range("F351").Value = 150574
$ws.Range("F352").Value = 307545
$ws.Range("F353").Value = 723593
$ws.Range("F354").Value = 316805
$ws.Range("F355").Value = 222145
$ws.Range("F356").Value = 160573
$ws.Range("F357").Value = 138374
$ws.Range("F358").Value = 159111
$ws.Range("F359").Value = 321217
$ws.Range("F360").Value = 749796
$ws.Range("F361").Value = 332926
$ws.Range("F362").Value = 229018
$ws.Range("F363").Value = 188786
$ws.Range("F364").Value = 168792
$ws.Range("F365").Value = 184709
$ws.Range("F366").Value = 339534
$ws.Range("F367").Value = 767148
$ws.Range("F368").Value = 346284
$ws.Range("F369").Value = 235136
$ws.Range("F370").Value = 181017
$ws.Range("F371").Value = 160269
$ws.Range("F372").Value = 178556
$ws.Range("F373").Value = 350267
$ws.Range("F374").Value = 773768
$ws.Range("F375").Value = 351322
$ws.Range("F376").Value = 222193
$ws.Range("F377").Value = 176656
$ws.Range("F378").Value = 157366
$ws.Range("F379").Value = 179772
$ws.Range("F380").Value = 344761
$ws.Range("F381").Value = 746792
$ws.Range("F382").Value = 356985
$ws.Range("F383").Value = 221106
$ws.Range("F384").Value = 172001
$ws.Range("F385").Value = 150871
$ws.Range("F386").Value = 182878
$ws.Range("F387").Value = 351605
$ws.Range("F388").Value = 730976
$ws.Range("F389").Value = 353696
$ws.Range("F390").Value = 219950
$ws.Range("F391").Value = 177665
$ws.Range("F392").Value = 221604
$ws.Range("F393").Value = 308176
$ws.Range("F394").Value = 166042
$ws.Range("F395").Value = 752521
$ws.Range("F396").Value = 164820
$ws.Range("F397").Value = 107847
$ws.Range("F398").Value = 298786
$ws.Range("F399").Value = 200509
$ws.Range("F400").Value = 150570
$ws.Range("F401").Value = 272570
$ws.Range("G401").Value = 932
$ws.Range("F402").Value = 720164
$ws.Range("F403").Value = 352221
$ws.Range("F404").Value = 225620
$ws.Range("F405").Value = 174097
$ws.Range("G405").Value = 693
$ws.Range("F406").Value = 170937
$ws.Range("F407").Value = 158369
$ws.Range("F408").Value = 304535
$ws.Range("F409").Value = 706921
$ws.Range("F410").Value = 364141
$ws.Range("G410").Value = 635
$ws.Range("F411").Value = 225258
$ws.Range("G411").Value = 827
$ws.Range("F412").Value = 176234
$ws.Range("G412").Value = 646
$ws.Range("F413").Value = 149560
$ws.Range("F414").Value = 148834
$ws.Range("F415").Value = 307538
$ws.Range("G415").Value = 694
$ws.Range("F416").Value = 669656
$ws.Range("G416").Value = 931
$ws.Range("F417").Value = 341551
$ws.Range("F418").Value = 202644
$ws.Range("F419").Value = 149307
$ws.Range("F420").Value = 138735
$ws.Range("G420").Value = 500
$ws.Range("F421").Value = 152880
$ws.Range("G421").Value = 531
$ws.Range("F422").Value = 297702
$ws.Range("G422").Value = 645
$ws.Range("F427").Value = 90347
$ws.Range("F430").Value = 173660
$ws.Range("F433").Value = 85767
$ws.Range("G433").Value = 264
$ws.Range("F434").Value = 80079
$ws.Range("F435").Value = 82277
$ws.Range("G435").Value = 266
$ws.Range("F436").Value = 144494
$ws.Range("G436").Value = 351
$ws.Range("F437").Value = 166541
$ws.Range("G437").Value = 275
$ws.Range("F438").Value = 121449
$ws.Range("G438").Value = 251
$ws.Range("F439").Value = 88919
$ws.Range("G439").Value = 318
$ws.Range("F440").Value = 73878
$ws.Range("G440").Value = 224
$ws.Range("F441").Value = 68230
$ws.Range("G441").Value = 202
$ws.Range("F442").Value = 69682
$ws.Range("G442").Value = 168

# --- New rows 443-445 ---
$ws.Range("A443").Value = 44337
$ws.Range("B443").Value = 388719
$ws.Range("C443").Value = 5671
$ws.Range("D443").Value = 190
$ws.Range("E443").Value = 12286
$ws.Range("F443").Value = 98429
$ws.Range("G443").Value = 297

$ws.Range("A444").Value = 44338
$ws.Range("B444").Value = 388835
$ws.Range("C444").Value = 4122
$ws.Range("D444").Value = 116
$ws.Range("E444").Value = 12292
$ws.Range("F444").Value = 89502
$ws.Range("G444").Value = 247

$ws.Range("A445").Value = 44339
$ws.Range("B445").Value = 388854
$ws.Range("C445").Value = 1167
$ws.Range("D445").Value = 19
$ws.Range("E445").Value = 12296
$ws.Range("F445").Value = 69004
$ws.Range("G445").Value = 324
